$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.295.56"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.167.47"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.10%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.42"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.69%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.167.24"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.35%  "

$ws.Range("E10").Value = "  -1.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.29"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("E13").Value = "  +0.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.76"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.689.17"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.11%  "

$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.165.90"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.265.06"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.57"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.98"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.696"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.64"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.63%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.22"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.95%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.07%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("E27").Value = "  -0.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.07"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("E30").Value = "  -2.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.72"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.16"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.01%  "

$ws.Range("E33").Value = "  -2.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.42"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.85%  "

$ws.Range("E35").Value = "  -2.90%  "

$ws.Range("E36").Value = "  -0.16%  "

$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0733"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.14%  "

$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.36"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0390"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.15"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.55%  "

$ws.Range("E41").Value = "  -2.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "392.54"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.790.28"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.78%  "

$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "126.98"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.38%  "

$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.73"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.11"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.37%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.06"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.54%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.112"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.41%  "
